$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,3).Value = 46064
$ws.Cells.Item(3,3).Value = 46064
$ws.Cells.Item(4,1).Value = "A 35197-2025"
$ws.Cells.Item(4,2).Value = 45853
$ws.Cells.Item(4,3).Value = 46064
$ws.Cells.Item(4,7).Value = 0.9
$ws.Cells.Item(4,8).Value = 1
$ws.Cells.Item(4,9).Value = 0
$ws.Cells.Item(4,10).Value = 3
$ws.Cells.Item(4,15).Value = 3
$ws.Cells.Item(4,18).Value = "Grönsångare`r`nTallticka`r`nVintertagging"
$ws.Cells.Item(4,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/artfynd/A 35197-2025 artfynd.xlsx", "A 35197-2025")'
$ws.Cells.Item(4,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/kartor/A 35197-2025 karta.png", "A 35197-2025")'
$ws.Cells.Item(4,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomål/A 35197-2025 FSC-klagomål.docx", "A 35197-2025")'
$ws.Cells.Item(4,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomålsmail/A 35197-2025 FSC-klagomål mail.docx", "A 35197-2025")'
$ws.Cells.Item(4,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsyn/A 35197-2025 tillsynsbegäran.docx", "A 35197-2025")'
$ws.Cells.Item(4,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsynsmail/A 35197-2025 tillsynsbegäran mail.docx", "A 35197-2025")'
$ws.Cells.Item(4,26).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/fåglar/A 35197-2025 prioriterade fågelarter.docx", "A 35197-2025")'
$ws.Cells.Item(5,1).Value = "A 37417-2023"
$ws.Cells.Item(5,2).Value = 45155
$ws.Cells.Item(5,3).Value = 46064
$ws.Cells.Item(5,7).Value = 12.9
$ws.Cells.Item(5,8).Value = 2
$ws.Cells.Item(5,10).Value = 1
$ws.Cells.Item(5,11).Value = 1
$ws.Cells.Item(5,15).Value = 2
$ws.Cells.Item(5,16).Value = 1
$ws.Cells.Item(5,18).Value = "Knärot`r`nUllticka`r`nBlåsippa"
$ws.Cells.Item(5,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/artfynd/A 37417-2023 artfynd.xlsx", "A 37417-2023")'
$ws.Cells.Item(5,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/kartor/A 37417-2023 karta.png", "A 37417-2023")'
$ws.Cells.Item(5,21).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/knärot/A 37417-2023 karta knärot.png", "A 37417-2023")'
$ws.Cells.Item(5,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomål/A 37417-2023 FSC-klagomål.docx", "A 37417-2023")'
$ws.Cells.Item(5,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomålsmail/A 37417-2023 FSC-klagomål mail.docx", "A 37417-2023")'
$ws.Cells.Item(5,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsyn/A 37417-2023 tillsynsbegäran.docx", "A 37417-2023")'
$ws.Cells.Item(5,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsynsmail/A 37417-2023 tillsynsbegäran mail.docx", "A 37417-2023")'
$ws.Cells.Item(5,26).ClearContents()
$ws.Cells.Item(6,1).Value = "A 47653-2024"
$ws.Cells.Item(6,2).Value = 45588
$ws.Cells.Item(6,3).Value = 46064
$ws.Cells.Item(6,7).Value = 3
$ws.Cells.Item(6,9).Value = 1
$ws.Cells.Item(6,11).Value = 0
$ws.Cells.Item(6,15).Value = 1
$ws.Cells.Item(6,16).Value = 0
$ws.Cells.Item(6,18).Value = "Backklöver`r`nNästrot`r`nBlåsippa"
$ws.Cells.Item(6,19).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/artfynd/A 47653-2024 artfynd.xlsx", "A 47653-2024")'
$ws.Cells.Item(6,20).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/kartor/A 47653-2024 karta.png", "A 47653-2024")'
$ws.Cells.Item(6,21).ClearContents()
$ws.Cells.Item(6,22).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomål/A 47653-2024 FSC-klagomål.docx", "A 47653-2024")'
$ws.Cells.Item(6,23).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/klagomålsmail/A 47653-2024 FSC-klagomål mail.docx", "A 47653-2024")'
$ws.Cells.Item(6,24).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsyn/A 47653-2024 tillsynsbegäran.docx", "A 47653-2024")'
$ws.Cells.Item(6,25).Formula = '=HYPERLINK("https://klasma.github.io/Logging_0305/tillsynsmail/A 47653-2024 tillsynsbegäran mail.docx", "A 47653-2024")'
$ws.Cells.Item(7,3).Value = 46064
$ws.Cells.Item(8,3).Value = 46064
$ws.Cells.Item(9,3).Value = 46064
$ws.Cells.Item(10,3).Value = 46064
$ws.Cells.Item(11,3).Value = 46064
$ws.Cells.Item(12,3).Value = 46064
$ws.Cells.Item(13,1).Value = "A 12146-2023"
$ws.Cells.Item(13,2).Value = 44998.47842592592
$ws.Cells.Item(13,3).Value = 46064
$ws.Cells.Item(13,7).Value = 3.1
$ws.Cells.Item(14,1).Value = "A 32023-2023"
$ws.Cells.Item(14,2).Value = 45119.49833333334
$ws.Cells.Item(14,3).Value = 46064
$ws.Cells.Item(14,7).Value = 3.1
$ws.Cells.Item(15,1).Value = "A 4524-2024"
$ws.Cells.Item(15,2).Value = 45327
$ws.Cells.Item(15,3).Value = 46064
$ws.Cells.Item(15,7).Value = 4.6
$ws.Cells.Item(16,1).Value = "A 23370-2025"
$ws.Cells.Item(16,2).Value = 45791.70907407408
$ws.Cells.Item(16,3).Value = 46064
$ws.Cells.Item(16,7).Value = 3.8
$ws.Cells.Item(17,1).Value = "A 37415-2023"
$ws.Cells.Item(17,2).Value = 45155
$ws.Cells.Item(17,3).Value = 46064
$ws.Cells.Item(17,7).Value = 6.6
$ws.Cells.Item(18,1).Value = "A 55962-2023"
$ws.Cells.Item(18,2).Value = 45240
$ws.Cells.Item(18,3).Value = 46064
$ws.Cells.Item(18,6).ClearContents()
$ws.Cells.Item(18,7).Value = 3.4
$ws.Cells.Item(19,1).Value = "A 11989-2025"
$ws.Cells.Item(19,2).Value = 45728.60074074074
$ws.Cells.Item(19,3).Value = 46064
$ws.Cells.Item(19,7).Value = 9.6
$ws.Cells.Item(20,3).Value = 46064
$ws.Cells.Item(21,1).Value = "A 35198-2025"
$ws.Cells.Item(21,2).Value = 45853
$ws.Cells.Item(21,3).Value = 46064
$ws.Cells.Item(21,6).ClearContents()
$ws.Cells.Item(21,7).Value = 1.2
$ws.Cells.Item(22,1).Value = "A 35300-2025"
$ws.Cells.Item(22,2).Value = 45854.41511574074
$ws.Cells.Item(22,3).Value = 46064
$ws.Cells.Item(22,7).Value = 2.2
$ws.Cells.Item(23,1).Value = "A 35193-2025"
$ws.Cells.Item(23,2).Value = 45853
$ws.Cells.Item(23,3).Value = 46064
$ws.Cells.Item(23,7).Value = 1.9
$ws.Cells.Item(24,1).Value = "A 37410-2023"
$ws.Cells.Item(24,2).Value = 45155
$ws.Cells.Item(24,3).Value = 46064
$ws.Cells.Item(24,7).Value = 20.9
$ws.Cells.Item(25,1).Value = "A 12156-2023"
$ws.Cells.Item(25,2).Value = 44998.49157407408
$ws.Cells.Item(25,3).Value = 46064
$ws.Cells.Item(25,7).Value = 0.5
$ws.Cells.Item(26,3).Value = 46064
$ws.Cells.Item(27,1).Value = "A 34202-2022"
$ws.Cells.Item(27,2).Value = 44791.64837962963
$ws.Cells.Item(27,3).Value = 46064
$ws.Cells.Item(27,7).Value = 2
$ws.Cells.Item(28,1).Value = "A 12154-2023"
$ws.Cells.Item(28,2).Value = 44998
$ws.Cells.Item(28,3).Value = 46064
$ws.Cells.Item(28,7).Value = 2.7
$ws.Cells.Item(29,1).Value = "A 67005-2021"
$ws.Cells.Item(29,2).Value = 44522
$ws.Cells.Item(29,3).Value = 46064
$ws.Cells.Item(29,7).Value = 1.3
$ws.Cells.Item(30,1).Value = "A 15732-2025"
$ws.Cells.Item(30,2).Value = 45747
$ws.Cells.Item(30,3).Value = 46064
$ws.Cells.Item(30,6).Value = "Kyrkan"
$ws.Cells.Item(30,7).Value = 1.4
$ws.Cells.Item(31,1).Value = "A 23250-2022"
$ws.Cells.Item(31,2).Value = 44719
$ws.Cells.Item(31,3).Value = 46064
$ws.Cells.Item(31,6).Value = "Naturvårdsverket"
$ws.Cells.Item(31,7).Value = 1
